$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the designation for the HOD (row 2, column B) from
# "HOD & PROFESSOR" to "Professor & Head"
$ws.Range("B2").Value = "Professor & Head"

# Enable word-wrap on the photo-path cell for row 6 (its text contains
# an embedded line break), matching the wrap formatting used elsewhere
$ws.Range("C6").WrapText = $true
